$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Uniform - Random")

# Fill in the previously-zeroed sample values for test #6 (row 7) and
# test #7 (row 8, Edges column) on the "Uniform - Random" sheet.
$ws.Range("D7").Value = 29388
$ws.Range("E7").Value = 75419
$ws.Range("C8").Value = 44822

# Move the active selection to D8, matching where the user left off editing.
$ws.Range("D8").Select()
